# Daily attendance processing - 2025-12-03 21:26:16
# Normalizes the "Recorded By" (column G) values so that when the list of
# recorders is exactly "<someone>, System", the order is flipped to
# "System, <someone>" (System is always listed first when paired with a
# single other recorder).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $value = $cell.Value()

    if ($null -ne $value -and $value -is [string]) {
        $parts = $value -split ', '
        if ($parts.Count -eq 2 -and $parts[1] -eq 'System') {
            $cell.Value = "System, " + $parts[0]
        }
    }
}
